$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = "dct:creator"
$ws.Cells.Item(13, 2).Value = "Minka"
